# The commit swaps the two theme parts of the deck:
#   ppt/theme/theme1.xml  (the "Integral" theme, used by the one SlideMaster
#                           and therefore by every slide) ends up holding the
#                           colours that used to live in theme2.xml ("Office
#                           Theme").
#   ppt/theme/theme2.xml  (only referenced by the NotesMaster) ends up
#                           holding the colours that used to live in
#                           theme1.xml ("Integral").
#
# The font scheme (majorFont/minorFont) and format scheme (fill/line/effect
# styles) are already byte-for-byte identical between the two theme parts,
# so the only real content difference is the <a:clrScheme> (12 colours) and
# the name= attributes on <a:theme>/<a:clrScheme>.
#
# PowerPoint's object model exposes the *active* theme's colour scheme via
# SlideMaster.Theme.ThemeColorScheme (a 12-entry collection of RGBColor
# objects, in the fixed order dk1, lt1, dk2, lt2, accent1-6, hlink,
# folHlink). Writing each .RGB value rewrites the corresponding
# <a:srgbClr val="..."/> inside ppt/theme/theme1.xml, which is the theme
# part that actually drives every slide's appearance.

$p = $ppt.ActivePresentation
$scheme = $p.SlideMaster.Theme.ThemeColorScheme

function Set-ThemeColor($scheme, $index, $r, $g, $b) {
    $packed = $r + ($g * 256) + ($b * 65536)
    $scheme.Item($index).RGB = $packed
}

# Target palette = the "Office Theme" colours (formerly theme2.xml),
# now written into theme1.xml per the diff.
Set-ThemeColor $scheme 1  0x00 0x00 0x00   # dk1
Set-ThemeColor $scheme 2  0xFF 0xFF 0xFF   # lt1
Set-ThemeColor $scheme 3  0x44 0x54 0x6A   # dk2
Set-ThemeColor $scheme 4  0xE7 0xE6 0xE6   # lt2
Set-ThemeColor $scheme 5  0x5B 0x9B 0xD5   # accent1
Set-ThemeColor $scheme 6  0xED 0x7D 0x31   # accent2
Set-ThemeColor $scheme 7  0xA5 0xA5 0xA5   # accent3
Set-ThemeColor $scheme 8  0xFF 0xC0 0x00   # accent4
Set-ThemeColor $scheme 9  0x44 0x72 0xC4   # accent5
Set-ThemeColor $scheme 10 0x70 0xAD 0x47   # accent6
Set-ThemeColor $scheme 11 0x05 0x63 0xC1   # hlink
Set-ThemeColor $scheme 12 0x95 0x4F 0x72   # folHlink
